# Applies the crypto-price refresh for cryptos.xlsx (Sheet1).
# Plain numeric-looking strings (e.g. "226.09") are written via a
# text-format round-trip so they stay text cells (matching the sheet's
# existing inlineStr/text convention) instead of being auto-parsed as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '34.418.21'
$ws.Range('E2').Value = '  +0.71%  '
# Row 3
$ws.Range('D3').Value = '1.790.35'
$ws.Range('E3').Value = '  +0.13%  '
# Row 4
$ws.Range('E4').Value = '  -0.02%  '
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '226.09'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.18%  '
# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.556'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +1.50%  '
# Row 7
$ws.Range('E7').Value = '  -0.05%  '
# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '32.64'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +2.48%  '
# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.296'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +1.34%  '
# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0692'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +0.43%  '
# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0949'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.44%  '
# Row 12
$ws.Range('D12').Value = '2.051.24'
# Row 13
$ws.Range('B13').Value = 'Chainlink'
$ws.Range('C13').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.03'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.12%  '
# Row 14
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.787.35'
$ws.Range('E14').Value = '  -0.22%  '
# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.634'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +1.80%  '
# Row 16
$ws.Range('D16').Value = '34.418.60'
# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.26'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +1.88%  '
# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '68.76'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.76%  '
# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '246.64'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.18%  '
# Row 20
$ws.Range('E20').Value = '  +2.87%  '
# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.16'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +3.12%  '
# Row 22
$ws.Range('E22').Value = '  -0.12%  '
# Row 23
$ws.Range('E23').Value = '  +1.33%  '
# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.08'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.74%  '
# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '164.80'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +2.33%  '
# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.22'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.78%  '
# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.48'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.90%  '
# Row 28
$ws.Range('E28').Value = '  +2.44%  '
# Row 29
$ws.Range('E29').Value = '  +0.03%  '
# Row 30
$ws.Range('E30').Value = '  +3.63%  '
# Row 31
$ws.Range('E31').Value = '  -0.18%  '
# Row 33
$ws.Range('E33').Value = '  +7.08%  '
# Row 34
$ws.Range('E34').Value = '  +1.24%  '
# Row 35
$ws.Range('D35').Value = '1.419.91'
$ws.Range('E35').Value = '  -1.83%  '
# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.57'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +5.99%  '
# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.665'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +2.58%  '
# Row 38
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0192'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +0.23%  '
# Row 39
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.06'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +1.74%  '
# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '84.63'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +5.29%  '
# Row 41
$ws.Range('E41').Value = '  +0.68%  '
# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.934'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +1.25%  '
# Row 43
$ws.Range('E43').Value = '  +1.78%  '
# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.56'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.42%  '
# Row 45
$ws.Range('E45').Value = '  +3.10%  '
# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '6.08'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.44%  '
# Row 47
$ws.Range('E47').Value = '  +0.08%  '
# Row 48
$ws.Range('D48').Value = '1.948.08'
$ws.Range('E48').Value = '  +0.01%  '
# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '105.43'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.14%  '
# Row 50
$ws.Range('E50').Value = '  -0.03%  '
# Row 51
$ws.Range('E51').Value = '  -4.42%  '
